$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 148
$ws.Range("J2").Value = 694
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 118
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 84
$ws.Range("T2").Value = 115
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 1089
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1091
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 2
